$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row with two new columns (P1 = 14, Q1 = 15),
# matching the existing header formatting (bold, bordered, centered).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the parallel/inner-loop columns (I, K, M, O) and append two
# new columns (P, Q) for every data row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O

    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
